# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, a new (blank) column is inserted
# immediately to the left of the existing "Late" column (column N),
# pushing "Late" / "heading" / "Outstanding" one column to the right
# (N->O, O->P, P->Q). The new column inherits the width of the column
# to its left ("In Advance", column M). The "Repayment schedule" sheet
# also becomes the active sheet/tab, with a new selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

$existingWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = $existingWidth

$ws.Activate() | Out-Null
$ws.Range("S10").Select() | Out-Null
